# Anjana's diary entry for week 4
# Adds four new diary entries (rows 33, 35, 37, 39, 41) to the log table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33 copies its formatting (date style, wrap-text body style, mood style)
# from the previous filled-in entry (row 32), then gets its own values.
$ws.Range("A32:G32").Copy()
$ws.Range("A33:G33").PasteSpecial(-4122)
$ws.Range("A33").Value = 43860
$ws.Range("B33").Value = "5pm - 7.50pm"
$ws.Range("C33").Value = "N/A"
$ws.Range("D33").Value = "Learn more concepts and practice the same"
$ws.Range("E33").Value = "Learned about key expert practices,structural vs behavioral models, UML,call graphs, sequence diagrams etc. "
$ws.Range("F33").Value = "Understood the importance of focusing on the core essence of the system rather than the extra functionalities, going deeper into the code base as when needed and how working with others can help us"
$ws.Range("G33").Value = "Happy to learn different ways to model code"
$ws.Rows.Item(33).RowHeight = 109.2

# The blank spacer row right after (row 34) keeps its date cell's style in
# sync with the date column even though there's no value in it.
$ws.Range("A33").Copy()
$ws.Range("A34").PasteSpecial(-4122)

# Row 35
$ws.Range("A32:G32").Copy()
$ws.Range("A35:G35").PasteSpecial(-4122)
$ws.Range("A35").Value = 43864
$ws.Range("B35").Value = "11:00 am - 1:00pm and 9:30 pm - 11:00 pm"
$ws.Range("C35").Value = "Aman, Vaishakhi"
$ws.Range("D35").Value = "Decide two features for the homework and start working on it."
$ws.Range("E35").Value = "We decided two features: `n1. how does the h2 database support embedded and server mode`n2. How is data actually persisted from h2 onto our disk.`nWe were also able to figure out the implementation of the first feature"
$ws.Range("F35").Value = "We were facing some issues with running our application because Tools.jar file was not being detected in pom.xml. So we had to downgrade our java version, edit the system path in pom.xml to reflect the location of tools.jar and do a maven clean.`nSince the code uses proper naming conventions it was much easier to read the codebase."
$ws.Range("G35").Value = "Was little frustrated initially when we were not able to run the system. Happy to have figured out the first feature"
$ws.Rows.Item(35).RowHeight = 30.6

# Row 37
$ws.Range("A32:G32").Copy()
$ws.Range("A37:G37").PasteSpecial(-4122)
$ws.Range("A37").Value = 43865
$ws.Range("B37").Value = "10pm - 12pm"
$ws.Range("C37").Value = "Aman, Vaishakhi"
$ws.Range("D37").Value = "Figure out second feature"
$ws.Range("E37").Value = "We were not able to figure out our second feature, even after spending a lot of time and doing a lot of debugging Hence we decided to choose a different feature as our second option which was ""How does H2 database handle multiple commands like Insert/Delete etc"""
$ws.Range("F37").Value = "Some features are very difficult to figure out "
$ws.Range("G37").Value = "Slightly nervous if we will be able to complete the homework or not"
$ws.Rows.Item(37).RowHeight = 124.8

# Row 39
$ws.Range("A32:G32").Copy()
$ws.Range("A39:G39").PasteSpecial(-4122)
$ws.Range("A39").Value = 43866
$ws.Range("B39").Value = "10pm - 2am"
$ws.Range("C39").Value = "Aman, Vaishakhi"
$ws.Range("D39").Value = "Finish working on second feature"
$ws.Range("E39").Value = "Able to understand the flow and we could draw the diagrams"
$ws.Range("F39").Value = "It was difficult and had to spend so much time"
$ws.Range("G39").Value = "Relaxed as we were able to finally finish the second feature"
$ws.Rows.Item(39).RowHeight = 31.2

# Row 41
$ws.Range("A32:G32").Copy()
$ws.Range("A41:G41").PasteSpecial(-4122)
$ws.Range("A41").Value = 43867
$ws.Range("B41").Value = "10am-12.30pm"
$ws.Range("C41").Value = "Aman, Vaishakhi"
$ws.Range("D41").Value = "Prepare reports"
$ws.Range("E41").Value = "Finished writing the reports"
$ws.Range("F41").Value = "Since we have prepared the flow digrams while looking at the code, it was easier to write the report"
$ws.Range("G41").Value = "Happy to finish the homework"
$ws.Rows.Item(41).RowHeight = 46.8

# Scroll position/selection left where the author was last working.
$ws.Application.Goto($ws.Range("A6"))
$ws.Range("G41").Select()
